$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.08076402653203729
$ws.Range("D2").Value = 0.175567214580866
$ws.Range("E2").Value = 0.1530259809589367
$ws.Range("F2").Value = 1.445659671885061
$ws.Range("G2").Value = 0.00246107558446808
$ws.Range("I2").Value = 0.5250896831123804
$ws.Range("J2").Value = 0.1729647364851488
$ws.Range("K2").Value = 0.614301067686938
$ws.Range("M2").Value = 0.2847130190385627
$ws.Range("N2").Value = 1.606286738031979
$ws.Range("O2").Value = 3.504915846646043
$ws.Range("B3").Value = 0.07103114751537021
$ws.Range("D3").Value = 0.1733145716070226
$ws.Range("E3").Value = 0.1524065783082875
$ws.Range("F3").Value = 1.44628900331282
$ws.Range("G3").Value = 0.002463692435612008
$ws.Range("I3").Value = 0.5315640392458052
$ws.Range("J3").Value = 0.1731919482371644
$ws.Range("K3").Value = 0.5495467292039109
$ws.Range("M3").Value = 0.2688932535344222
$ws.Range("N3").Value = 1.621621793670325
$ws.Range("O3").Value = 3.512950023484478
$ws.Range("B4").Value = 0.06505023982553837
$ws.Range("D4").Value = 0.171997797137756
$ws.Range("E4").Value = 0.1520919391578275
$ws.Range("F4").Value = 1.447414254614174
$ws.Range("G4").Value = 0.002465386315214413
$ws.Range("I4").Value = 0.5357740824352826
$ws.Range("J4").Value = 0.1734112798005754
$ws.Range("K4").Value = 0.5098179786788535
$ws.Range("M4").Value = 0.2592774678017378
$ws.Range("N4").Value = 1.631559231558565
$ws.Range("O4").Value = 3.519745360594442
$ws.Range("B5").Value = 0.06261194150879135
$ws.Range("D5").Value = 0.1714779607190806
$ws.Range("E5").Value = 0.1519802700391217
$ws.Range("F5").Value = 1.448058711018312
$ws.Range("G5").Value = 0.002466098557368284
$ws.Range("I5").Value = 0.5375487565399615
$ws.Range("J5").Value = 0.1735207480988983
$ws.Range("K5").Value = 0.4936367820641863
$ws.Range("M5").Value = 0.2553837403482007
$ws.Range("N5").Value = 1.635740056022637
$ws.Range("O5").Value = 3.522982911517659
$ws.Range("B6").Value = 0.06220700753920028
$ws.Range("D6").Value = 0.1713926567025368
$ws.Range("E6").Value = 0.1519627279793738
$ws.Range("F6").Value = 1.448176954416624
$ws.Range("G6").Value = 0.002466218153562182
$ws.Range("I6").Value = 0.5378470065614511
$ws.Range("J6").Value = 0.1735401390444373
$ws.Range("K6").Value = 0.4909504514899652
$ws.Range("M6").Value = 0.2547386937132856
$ws.Range("N6").Value = 1.636442207114925
$ws.Range("O6").Value = 3.523548797357364
$ws.Range("B7").Value = 0.065017359990307
$ws.Range("D7").Value = 0.1719907184831229
$ws.Range("E7").Value = 0.1520903660978696
$ws.Range("F7").Value = 1.447422193088336
$ws.Range("G7").Value = 0.002465395831572625
$ws.Range("I7").Value = 0.5357977772109743
$ws.Range("J7").Value = 0.1734126747667055
$ws.Range("K7").Value = 0.5095997174508966
$ws.Range("M7").Value = 0.2592248549702063
$ws.Range("N7").Value = 1.631615084156103
$ws.Range("O7").Value = 3.519787126692279
$ws.Range("B8").Value = 0.07740925820723987
$ws.Range("D8").Value = 0.1747767740065171
$ws.Range("E8").Value = 0.152798799604227
$ws.Range("F8").Value = 1.445723387667904
$ws.Range("G8").Value = 0.002461959831536908
$ws.Range("I8").Value = 0.5272733254996433
$ws.Range("J8").Value = 0.1730265216665785
$ws.Range("K8").Value = 0.5919679120618468
$ws.Range("M8").Value = 0.2792382411690539
$ws.Range("N8").Value = 1.611466035990393
$ws.Range("O8").Value = 3.507299558937945
$ws.Range("B9").Value = 0.1016634899278444
$ws.Range("D9").Value = 0.1807639182862886
$ws.Range("E9").Value = 0.1547078849242247
$ws.Range("F9").Value = 1.448250828886557
$ws.Range("G9").Value = 0.002455910159139662
$ws.Range("I9").Value = 0.5124187786961985
$ws.Range("J9").Value = 0.1729020572964259
$ws.Range("K9").Value = 0.7537044064186773
$ws.Range("M9").Value = 0.3192506650167601
$ws.Range("N9").Value = 1.576089563627839
$ws.Range("O9").Value = 3.497587355974929
$ws.Range("B10").Value = 0.1194464784621658
$ws.Range("D10").Value = 0.185478766887158
$ws.Range("E10").Value = 0.1564260406315192
$ws.Range("F10").Value = 1.453677170223429
$ws.Range("G10").Value = 0.002451880930335281
$ws.Range("I10").Value = 0.5026387123383049
$ws.Range("J10").Value = 0.1731958426696707
$ws.Range("K10").Value = 0.8726323179833457
$ws.Range("M10").Value = 0.3491069231723145
$ws.Range("N10").Value = 1.552614528228546
$ws.Range("O10").Value = 3.499463868327439
$ws.Range("B11").Value = 0.1275267931052753
$ws.Range("D11").Value = 0.1876916767888162
$ws.Range("E11").Value = 0.1572759254398797
$ws.Range("F11").Value = 1.456920525714722
$ws.Range("G11").Value = 0.002450137262173684
$ws.Range("I11").Value = 0.4984353273999371
$ws.Range("J11").Value = 0.173413047153673
$ws.Range("K11").Value = 0.9267518652663682
$ws.Range("M11").Value = 0.3627874982420849
$ws.Range("N11").Value = 1.542480379899612
$ws.Range("O11").Value = 3.502275745830872
$ws.Range("B12").Value = 0.1305850873758061
$ws.Range("D12").Value = 0.1885393747744644
$ws.Range("E12").Value = 0.1576075451568322
$ws.Range("F12").Value = 1.45826005628625
$ws.Range("G12").Value = 0.00244948974761976
$ws.Range("I12").Value = 0.4968789198265995
$ws.Range("J12").Value = 0.1735073004620205
$ws.Range("K12").Value = 0.9472474712938492
$ws.Range("M12").Value = 0.3679819866009737
$ws.Range("N12").Value = 1.538721123113369
$ws.Range("O12").Value = 3.503622142049181
$ws.Range("B13").Value = 0.1299265011350883
$ws.Range("D13").Value = 0.1883563768867873
$ws.Range("E13").Value = 0.1575356901706293
$ws.Range("F13").Value = 1.457966613796202
$ws.Range("G13").Value = 0.002449628634238498
$ws.Range("I13").Value = 0.4972125489807677
$ws.Range("J13").Value = 0.1734864677127348
$ws.Range("K13").Value = 0.9428333153493043
$ws.Range("M13").Value = 0.3668626442020297
$ws.Range("N13").Value = 1.539527262933618
$ws.Range("O13").Value = 3.503319647231905
$ws.Range("B14").Value = 0.127778433017923
$ws.Range("D14").Value = 0.1877612231036068
$ws.Range("E14").Value = 0.1573030120400318
$ws.Range("F14").Value = 1.457028498726672
$ws.Range("G14").Value = 0.002450083735097014
$ws.Range("I14").Value = 0.4983065728629867
$ws.Range("J14").Value = 0.1734205609539075
$ws.Range("K14").Value = 0.9284380224516724
$ws.Range("M14").Value = 0.3632145737943446
$ws.Range("N14").Value = 1.542169533880234
$ws.Range("O14").Value = 3.502380870954596
$ws.Range("B15").Value = 0.1264624716378364
$ws.Range("D15").Value = 0.1873979374181118
$ws.Range("E15").Value = 0.1571617635198805
$ws.Range("F15").Value = 1.456468373667221
$ws.Range("G15").Value = 0.002450364159447316
$ws.Range("I15").Value = 0.4989812943014975
$ws.Range("J15").Value = 0.173381753874331
$ws.Range("K15").Value = 0.9196206834536724
$ws.Range("M15").Value = 0.3609818345752416
$ws.Range("N15").Value = 1.543798201108839
$ws.Range("O15").Value = 3.501842516095053
$ws.Range("B16").Value = 0.1189182062995684
$ws.Range("D16").Value = 0.1853355119120863
$ws.Range("E16").Value = 0.1563718699676642
$ws.Range("F16").Value = 1.453480794701363
$ws.Range("G16").Value = 0.002451996674110504
$ws.Range("I16").Value = 0.5029183561747494
$ws.Range("J16").Value = 0.1731833281636881
$ws.Range("K16").Value = 0.8690957702832236
$ws.Range("M16").Value = 0.3482148317673222
$ws.Range("N16").Value = 1.553287779674484
$ws.Range("O16").Value = 3.499319510911789
$ws.Range("B17").Value = 0.1142875204778733
$ws.Range("D17").Value = 0.1840876691276492
$ws.Range("E17").Value = 0.1559047624936944
$ws.Range("F17").Value = 1.451846417897812
$ws.Range("G17").Value = 0.00245302098504726
$ws.Range("I17").Value = 0.5053965269611158
$ws.Range("J17").Value = 0.1730829927977737
$ws.Range("K17").Value = 0.8381044554315338
$ws.Range("M17").Value = 0.3404078142274685
$ws.Range("N17").Value = 1.55924884691197
$ws.Range("O17").Value = 3.498273286419476
$ws.Range("B18").Value = 0.111623214056948
$ws.Range("D18").Value = 0.1833763575491218
$ws.Range("E18").Value = 0.1556425241054207
$ws.Range("F18").Value = 1.450979308600324
$ws.Range("G18").Value = 0.002453618545891806
$ws.Range("I18").Value = 0.5068450235872834
$ws.Range("J18").Value = 0.1730331484029861
$ws.Range("K18").Value = 0.8202808838210558
$ws.Range("M18").Value = 0.3359267404412591
$ws.Range("N18").Value = 1.562728771271704
$ws.Range("O18").Value = 3.497855844437254
$ws.Range("B19").Value = 0.1107209847888839
$ws.Range("D19").Value = 0.1831366236938266
$ws.Range("E19").Value = 0.1555548399612121
$ws.Range("F19").Value = 1.450698250580658
$ws.Range("G19").Value = 0.002453822314907109
$ws.Range("I19").Value = 0.5073394301410765
$ws.Range("J19").Value = 0.1730176233155021
$ws.Range("K19").Value = 0.8142464697560285
$ws.Range("M19").Value = 0.3344111337728677
$ws.Range("N19").Value = 1.563915820977957
$ws.Range("O19").Value = 3.497746162478165
$ws.Range("B20").Value = 0.1147805552858756
$ws.Range("D20").Value = 0.1842198407861559
$ws.Range("E20").Value = 0.1559538216929184
$ws.Range("F20").Value = 1.452012851593835
$ws.Range("G20").Value = 0.00245291107623758
$ws.Range("I20").Value = 0.5051303287627427
$ws.Range("J20").Value = 0.1730928596775598
$ws.Range("K20").Value = 0.8414033532126552
$ws.Range("M20").Value = 0.3412379222456323
$ws.Range("N20").Value = 1.558608974326397
$ws.Range("O20").Value = 3.498365582460309
$ws.Range("B21").Value = 0.12840941634731
$ws.Range("D21").Value = 0.1879357710844261
$ws.Range("E21").Value = 0.1573710898984757
$ws.Range("F21").Value = 1.457301024748062
$ws.Range("G21").Value = 0.002449949714692401
$ws.Range("I21").Value = 0.4979842728405508
$ws.Range("J21").Value = 0.1734395937219446
$ws.Range("K21").Value = 0.9326662303216153
$ws.Range("M21").Value = 0.3642857239196502
$ws.Range("N21").Value = 1.541391309202698
$ws.Range("O21").Value = 3.502648969731325
$ws.Range("B22").Value = 0.1373075570401454
$ws.Range("D22").Value = 0.190420949405862
$ws.Range("E22").Value = 0.1583543862111227
$ws.Range("F22").Value = 1.461406099586824
$ws.Range("G22").Value = 0.00244808872733229
$ws.Range("I22").Value = 0.4935198039690594
$ws.Range("J22").Value = 0.1737361591831998
$ws.Range("K22").Value = 0.9923213003371814
$ws.Range("M22").Value = 0.379429994836336
$ws.Range("N22").Value = 1.530595106054371
$ws.Range("O22").Value = 3.507089789593778
$ws.Range("B23").Value = 0.1325593560809182
$ws.Range("D23").Value = 0.1890894088627277
$ws.Range("E23").Value = 0.1578243751057755
$ws.Range("F23").Value = 1.459155793373199
$ws.Range("G23").Value = 0.002449075180055786
$ws.Range("I23").Value = 0.4958837346377791
$ws.Range("J23").Value = 0.1735714801011952
$ws.Range("K23").Value = 0.960481718619036
$ws.Range("M23").Value = 0.3713398678353812
$ws.Range("N23").Value = 1.53631547267484
$ws.Range("O23").Value = 3.504569452155522
$ws.Range("B24").Value = 0.1145576606607932
$ws.Range("D24").Value = 0.1841600670055499
$ws.Range("E24").Value = 0.1559316223783682
$ws.Range("F24").Value = 1.451937381030675
$ws.Range("G24").Value = 0.002452960739059376
$ws.Range("I24").Value = 0.5052506029481343
$ws.Range("J24").Value = 0.1730883744401979
$ws.Range("K24").Value = 0.8399119409263278
$ws.Range("M24").Value = 0.3408626077006787
$ws.Range("N24").Value = 1.558898096166086
$ws.Range("O24").Value = 3.498323282122186
$ws.Range("B25").Value = 0.09510789804659225
$ws.Range("D25").Value = 0.1790884990091541
$ws.Range("E25").Value = 0.1541359141790899
$ws.Range("F25").Value = 1.44694009744255
$ws.Range("G25").Value = 0.002457473497759075
$ws.Range("I25").Value = 0.5162381543708143
$ws.Range("J25").Value = 0.1728680395950661
$ws.Range("K25").Value = 0.7099306838808275
$ws.Range("M25").Value = 0.3083450386779063
$ws.Range("N25").Value = 1.585217625616206
$ws.Range("O25").Value = 3.498632462916873
